# One-click update from Van Paper 08:11 AM on 2025-09-16
# Adds a new customer row ("TAQUERIA Y MERCADO ANDALE", customer #0008300)
# into the leaderboard sheet, right after "Shakopee Brewhall" (row 38),
# pushing the existing rows 39-42 down to 40-43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 39; this shifts the old rows
# 39-42 down to 40-43 and extends the sheet dimension to A1:F43.
$ws.Rows("39:39").Insert()

# Match the row height used throughout the rest of the table.
$ws.Rows("39:39").RowHeight = 13.05

# Populate the new row with the new customer's data.
$ws.Range("A39").Value = "TAQUERIA Y MERCADO ANDALE"
$ws.Range("B39").Value = "Pietrs, Josh"
$ws.Range("C39").Value = "030"
$ws.Range("E39").Value = "0008300"
# Column D (Last Invoice Date) is intentionally left blank for this
# new prospect, same as row 38 above it.
